# "third teams time variable" — reorders the list of third-place teams
# stored (as a stringified Python list) in column O for a handful of rows.
# The set of teams per row is unchanged; only the textual ordering differs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 40-43 and 53-62: ['Colombia', 'Argentina'] -> ['Argentina', 'Colombia']
$rowsPair = @(40,41,42,43,53,54,55,56,57,58,59,60,61,62)
foreach ($r in $rowsPair) {
    $ws.Range("O$r").Value = "['Argentina', 'Colombia']"
}

# Rows 63-73: ['Colombia', 'Argentina', 'Scotland', 'Austria'] -> ['Scotland', 'Austria', 'Argentina', 'Colombia']
$rowsQuad = @(63,64,65,66,67,68,69,70,71,72,73)
foreach ($r in $rowsQuad) {
    $ws.Range("O$r").Value = "['Scotland', 'Austria', 'Argentina', 'Colombia']"
}

# Row 82: ['South Korea', 'United States', 'Saudi Arabia', 'Cameroon'] -> ['United States', 'Saudi Arabia', 'Cameroon', 'South Korea']
$ws.Range("O82").Value = "['United States', 'Saudi Arabia', 'Cameroon', 'South Korea']"

# Row 102: ['Netherlands', 'Italy'] -> ['Italy', 'Netherlands']
$ws.Range("O102").Value = "['Italy', 'Netherlands']"
